$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D31").Value = "2016-03-08 08:14:01"
$wsZhCn.Range("D32").Value = "2016-03-08 08:14:01"
$wsZhCn.Range("G31").Value = "2016-03-08 08:14:21"
$wsZhCn.Range("G32").Value = "2016-03-08 08:14:21"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D31").Value = "2016-03-08 08:14:06"
$wsDeDe.Range("D32").Value = "2016-03-08 08:14:06"
$wsDeDe.Range("G31").Value = "2016-03-08 08:14:28"
$wsDeDe.Range("G32").Value = "2016-03-08 08:14:28"
